# Introduce a new "PL_retirement_adjustment" sheet, aligned with the other
# PL_* time-series-factor sheets, inserted right after "Info" and before
# "PL_cohabitation_adjustment".

$wb = $excel.ActiveWorkbook

# Insert the new worksheet immediately before "PL_cohabitation_adjustment"
# (i.e. right after "Info"), matching the target sheet order.
$beforeSheet = $wb.Worksheets.Item("PL_cohabitation_adjustment")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "PL_retirement_adjustment"

# Header row
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Value"

# Body rows: years 2010-2070 with placeholder value 0 in column B.
$startYear = 2010
$endYear = 2070
$row = 2
for ($year = $startYear; $year -le $endYear; $year++) {
    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = 0
    $row++
}
